# Atualização dos dados e melhorias no codigo
#
# - Insert two header rows above the existing municipality data:
#     row 1: machine-style column headers (Unnamed: 0/1/2), bold,
#            centered + top-aligned, thin border all around
#     row 2: human-readable headers (municipio / nº de casos / nº de óbitos)
# - Append a new "(vazio)" row at the bottom of the data with
#   4 casos and 1 óbito.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down by two rows so it now starts on row 3.
$ws.Rows("1:2").Insert()

# Row 1 - "Unnamed" placeholder headers
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "Unnamed: 1"
$ws.Range("C1").Value = "Unnamed: 2"

# Row 2 - descriptive headers
$ws.Range("A2").Value = "municipio "
$ws.Range("B2").Value = "nº de casos"
$ws.Range("C2").Value = "nº de óbitos"

# New trailing row summarising the blank/unassigned municipality
$ws.Range("A52").Value = "(vazio)"
$ws.Range("B52").Value = 4
$ws.Range("C52").Value = 1

# Formatting for the new header row (row 1): bold, thin border all
# around, centered horizontally, aligned to top vertically.
$hdr = $ws.Range("A1:C1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
